# Weekly update: insert a new price-report row for the current week.
# This pushes the existing data rows (old rows 4-59) down by one (to 5-60)
# and inserts a brand-new row 4 with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(4).EntireRow.Insert()

$ws.Cells.Item(4, 1).Value = 4
$ws.Cells.Item(4, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(4, 3).Value = "Los Lagos"
$ws.Cells.Item(4, 4).Value = 44691
$ws.Cells.Item(4, 5).Value = 10
$ws.Cells.Item(4, 6).Value = 100112043
$ws.Cells.Item(4, 7).Value = "Pepino dulce"
$ws.Cells.Item(4, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 70
$ws.Cells.Item(4, 11).Value = 17000
$ws.Cells.Item(4, 12).Value = 17000
$ws.Cells.Item(4, 13).Value = 17000
$ws.Cells.Item(4, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(4, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(4, 16).Value = 944
$ws.Cells.Item(4, 17).Value = 18
$ws.Cells.Item(4, 18).Value = "Hortaliza"
